$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -4383343486.795609
$ws.Range("C2").Value = -4385389245.879074

$ws.Range("B3").Value = -4379616096.379884
$ws.Range("C3").Value = -4383707625.275123

$ws.Range("B4").Value = -4372161263.97213
$ws.Range("C4").Value = -4380344364.676104

$ws.Range("B5").Value = -4363215374.307608
$ws.Range("C5").Value = -4376308417.828602
